$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.954.38'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '2.466.12'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '2.465.82'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.330'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.52%  '
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '68.837.99'
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '2.494.73'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').Value = '2.600.85'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '428.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('E35').Value = '  -3.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '156.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.299'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('E43').Value = '  -4.70%  '
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '132.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0714'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0912'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.28%  '
